# Excel XlHAlign / XlVAlign constants (standard VBA values):
#   xlCenter = -4108   xlLeft = -4131   xlRight = -4152
$xlCenter = -4108
$xlLeft   = -4131

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Unmerge the two columns (Date / Who) that will grow from row 3:5 to 3:6 ---
$ws.Range("A3:A5").UnMerge()
$ws.Range("B3:B5").UnMerge()

# --- New row 6: "Profiles, Personas and scenarios" project / "Personas" task ---
# (written in this order so the shared-string table is appended in the same
#  sequence as the source: time, project, task, feeling)
$ws.Range("E6").Value = "1hr 7min (20:36 - 21:43)"
$ws.Range("C6").Value = "Profiles, Personas and scenarios"
$ws.Range("D6").Value = "Personas"
$ws.Range("F6").Value = "Pleased, identifying the profile of possible users was a really interesting and pleasing task"

# --- Re-merge Date/Who down through the new row ---
$ws.Range("A3:A6").Merge()
$ws.Range("B3:B6").Merge()

# --- Realign: Date + Who columns become centered horizontally (was right/left) ---
$ws.Range("A3:A6").HorizontalAlignment = $xlCenter
$ws.Range("A3:A6").VerticalAlignment = $xlCenter

$ws.Range("B3:B6").HorizontalAlignment = $xlCenter
$ws.Range("B3:B6").VerticalAlignment = $xlCenter

# --- Project column (C3:C5) keeps its left/center look ---
$ws.Range("C3:C5").HorizontalAlignment = $xlLeft
$ws.Range("C3:C5").VerticalAlignment = $xlCenter

# --- Time cell for the new row mirrors the other "time" cells (E2/E4 style: h:mm) ---
$ws.Range("E6").NumberFormat = "h:mm"
